$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.482.97'
$ws.Range("E2").Value = '  +5.13%  '
$ws.Range("D3").Value = '3.634.49'
$ws.Range("E3").Value = '  +5.46%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.99'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '191.69'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.86%  '
$ws.Range("E7").Value = '  +2.09%  '
$ws.Range("D8").Value = '3.628.84'
$ws.Range("E8").Value = '  +5.44%  '
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("E10").Value = '  +2.94%  '
$ws.Range("E11").Value = '  +3.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '58.30'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.49%  '
$ws.Range("E13").Value = '  +3.91%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.92'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.99%  '
$ws.Range("D15").Value = '4.211.23'
$ws.Range("E15").Value = '  +5.44%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.76'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.78%  '
$ws.Range("D17").Value = '3.629.78'
$ws.Range("E17").Value = '  +5.30%  '
$ws.Range("D18").Value = '70.427.16'
$ws.Range("E18").Value = '  +5.36%  '
$ws.Range("E19").Value = '  +5.03%  '
$ws.Range("E20").Value = '  +0.41%  '
$ws.Range("E21").Value = '  +4.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '488.90'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.26%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '19.60'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +16.39%  '
$ws.Range("E24").Value = '  -0.32%  '
$ws.Range("E25").Value = '  +1.19%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '91.09'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.42%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.14'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.37'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.71'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.95%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.13'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.89'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +10.52%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '631.72'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +7.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.32'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.12%  '
$ws.Range("E34").Value = '  +7.27%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '66.13'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.96%  '
$ws.Range("B36").Value = 'TheGraph'
$ws.Range("C36").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.414'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.51%  '
$ws.Range("B37").Value = 'InjectiveProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '38.90'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.66%  '
$ws.Range("E38").Value = '  +6.95%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.147'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.95%  '
$ws.Range("B40").Value = 'Dai'
$ws.Range("C40").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("E41").Value = '  +0.36%  '
$ws.Range("D42").Value = '3.311.88'
$ws.Range("E42").Value = '  +3.63%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.12'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +7.26%  '
$ws.Range("E44").Value = '  +10.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0451'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.10%  '
$ws.Range("E46").Value = '  +2.93%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.28'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.14'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.92%  '
$ws.Range("E49").Value = '  -1.66%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.31'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.71%  '
$ws.Range("B51").Value = 'FirstDigitalUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.06%  '
